$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("G3").Value = 78
$wsExhibit.Range("F5").Value = 504
$wsExhibit.Range("G5").Value = 73
$wsExhibit.Range("F8").Value = 21
$wsExhibit.Range("F9").Value = 1023
$wsExhibit.Range("F10").Value = 815
$wsExhibit.Range("F11").Value = 238
$wsExhibit.Range("F14").Value = 821
$wsExhibit.Range("F15").Value = 278
$wsExhibit.Range("F16").Value = 580
$wsExhibit.Range("F18").Value = 1325
$wsExhibit.Range("F19").Value = 119
$wsExhibit.Range("F20").Value = 853
$wsExhibit.Range("F21").Value = 1177
$wsExhibit.Range("F22").Value = 2858
$wsExhibit.Range("F23").Value = 1404
$wsExhibit.Range("F24").Value = 698
$wsExhibit.Range("F25").Value = 190
$wsExhibit.Range("F26").Value = 1270
$wsExhibit.Range("F28").Value = 1011
$wsExhibit.Range("F29").Value = 356
$wsExhibit.Range("F30").Value = 3062
$wsExhibit.Range("F31").Value = 592
$wsExhibit.Range("F32").Value = 534
$wsExhibit.Range("F33").Value = 1389
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 77
$wsShow.Range("F10").Value = 154
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 735
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 735
$wsAll.Range("G4").Value = 78
$wsAll.Range("F5").Value = 77
$wsAll.Range("F7").Value = 504
$wsAll.Range("G7").Value = 73
$wsAll.Range("F15").Value = 21
$wsAll.Range("F16").Value = 1023
$wsAll.Range("F17").Value = 815
$wsAll.Range("F18").Value = 238
$wsAll.Range("F22").Value = 154
$wsAll.Range("F26").Value = 821
$wsAll.Range("F27").Value = 278
$wsAll.Range("F28").Value = 580
$wsAll.Range("F30").Value = 1325
$wsAll.Range("F31").Value = 119
$wsAll.Range("F32").Value = 853
$wsAll.Range("F33").Value = 1177
$wsAll.Range("F34").Value = 2858
$wsAll.Range("F35").Value = 1404
$wsAll.Range("F36").Value = 698
$wsAll.Range("F37").Value = 190
$wsAll.Range("F38").Value = 1270
$wsAll.Range("F42").Value = 1011
$wsAll.Range("F43").Value = 356
$wsAll.Range("F44").Value = 3062
$wsAll.Range("F45").Value = 592
$wsAll.Range("F46").Value = 534
$wsAll.Range("F47").Value = 1389
